$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 3, shifting the header + data block
# down by one row (matches the OOXML diff: old row 3 -> new row 4, etc.)
$ws.Rows.Item(3).Insert()

# Scroll / reselect as seen in the saved view state
$ws.Range("M6").Select()
